# Add "limitTime" timer field (new column D) to the Stage sheet.
$wb = $excel.ActiveWorkbook

$stage = $wb.Worksheets.Item("Stage")
$character = $wb.Worksheets.Item("Character")
$gimmick = $wb.Worksheets.Item("Gimmick")

# --- Stage sheet: add new column D "limitTime" ---

# Copy formatting from column C into column D so the new column matches the
# look of the rest of the table (header / blank / comment / type / value rows).
$stage.Range("C1:C5").Copy()
$stage.Range("D1:D5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$stage.Range("D1").Value = "limitTime"
$stage.Range("D4").Value = "int"
$stage.Range("D5").Value = 180

$stage.Columns.Item(4).ColumnWidth = 15

# Stage page setup (A4 portrait)
$stage.PageSetup.PaperSize = 9
$stage.PageSetup.Orientation = 1

# --- Sheet selections / active sheet ---
$character.Range("I7").Select()
$gimmick.Range("H21").Select()

$stage.Activate()
$stage.Range("D3").Select()
